$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.855.95"
$ws.Range("E2").Value = "  -2.01%  "
$ws.Range("D3").Value = "1.826.94"
$ws.Range("E3").Value = "  -2.29%  "
$ws.Range("D4").Value = "0.9999"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "239.54"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6880"
$ws.Range("E6").Value = "  -2.58%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "0.07627"
$ws.Range("E8").Value = "  -3.24%  "
$ws.Range("D9").Value = "0.3021"
$ws.Range("E9").Value = "  -4.13%  "
$ws.Range("D10").Value = "23.47"
$ws.Range("E10").Value = "  -4.77%  "
$ws.Range("D11").Value = "0.07731"
$ws.Range("E11").Value = "  -3.38%  "
$ws.Range("D12").Value = "1.828.23"
$ws.Range("E12").Value = "  -2.79%  "
$ws.Range("D13").Value = "5.041"
$ws.Range("E13").Value = "  -3.34%  "
$ws.Range("D14").Value = "90.11"
$ws.Range("E14").Value = "  -4.27%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6720"
$ws.Range("E15").Value = "  -4.63%  "
$ws.Range("D16").Value = "6.414"
$ws.Range("E16").Value = "  -1.20%  "
$ws.Range("D17").Value = "0.000008275"
$ws.Range("E17").Value = "  -1.19%  "
$ws.Range("D18").Value = "28.850.40"
$ws.Range("E18").Value = "  -2.30%  "
$ws.Range("D19").Value = "242.52"
$ws.Range("E19").Value = "  -5.34%  "
$ws.Range("D20").Value = "2.084.86"
$ws.Range("E20").Value = "  -2.77%  "
$ws.Range("D21").Value = "12.61"
$ws.Range("E21").Value = "  -4.39%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.000"
$ws.Range("E22").Value = "  -0.07%  "
$ws.Range("D23").Value = "7.396"
$ws.Range("E23").Value = "  -2.95%  "
$ws.Range("E24").Value = "  -0.08%  "
$ws.Range("D25").Value = "0.1467"
$ws.Range("E25").Value = "  -5.70%  "
$ws.Range("D26").Value = "160.31"
$ws.Range("E26").Value = "  -0.55%  "
$ws.Range("D27").Value = "8.702"
$ws.Range("E27").Value = "  -4.02%  "
$ws.Range("D28").Value = "18.13"
$ws.Range("E28").Value = "  -3.70%  "
$ws.Range("D29").Value = "1.528"
$ws.Range("E29").Value = "  +1.92%  "
$ws.Range("D30").Value = "4.194"
$ws.Range("E30").Value = "  -3.27%  "
$ws.Range("D31").Value = "4.146"
$ws.Range("E31").Value = "  -2.44%  "
$ws.Range("E32").Value = "  -1.62%  "
$ws.Range("E33").Value = "  -4.46%  "
$ws.Range("D34").Value = "0.7588"
$ws.Range("E34").Value = "  +1.41%  "
$ws.Range("D35").Value = "1.809"
$ws.Range("E35").Value = "  -4.78%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.140"
$ws.Range("E36").Value = "  -2.58%  "
$ws.Range("D37").Value = "2.686"
$ws.Range("E37").Value = "  -1.02%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01830"
$ws.Range("E38").Value = "  -2.49%  "
$ws.Range("D39").Value = "1.210.59"
$ws.Range("E39").Value = "  -4.34%  "
$ws.Range("D40").Value = "2.676"
$ws.Range("E40").Value = "  -2.65%  "
$ws.Range("D41").Value = "0.9153"
$ws.Range("E41").Value = "  +1.61%  "
$ws.Range("D42").Value = "108.51"
$ws.Range("E42").Value = "  -0.64%  "
$ws.Range("D43").Value = "0.9997"
$ws.Range("E43").Value = "  -0.06%  "
$ws.Range("D44").Value = "1.986.11"
$ws.Range("E44").Value = "  -2.59%  "
$ws.Range("E45").Value = "  -5.85%  "
$ws.Range("D46").Value = "0.5159"
$ws.Range("D47").Value = "9.418"
$ws.Range("E47").Value = "  -1.05%  "
$ws.Range("D48").Value = "5.286"
$ws.Range("E48").Value = "  -11.55%  "
$ws.Range("D49").Value = "62.56"
$ws.Range("E49").Value = "  -12.89%  "
$ws.Range("D50").Value = "1.721"
$ws.Range("E50").Value = "  -5.06%  "
$ws.Range("D51").Value = "0.4161"
$ws.Range("E51").Value = "  -3.74%  "
